# Translate the dataset's Portuguese-facing labels (header + emotion
# categories) from English to Portuguese. The free-text review content
# in column B is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Emoção"
$ws.Range("B1").Value = "Texto"

# Emotion category labels for each data row
$labels = @{
    2  = "Positiva"
    3  = "Positiva"
    4  = "Positiva"
    5  = "Positiva"
    6  = "Positiva"
    7  = "Negativa"
    8  = "Negativa"
    9  = "Negativa"
    10 = "Negativa"
    11 = "Negativa"
    12 = "Neutra"
    13 = "Neutra"
    14 = "Neutra"
    15 = "Neutra"
    16 = "Neutra"
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
}

# Reflect the updated active cell selection in the saved sheet view
$ws.Range("A2").Select()
